# ---------------------------------------------------------------------------
# Applies the two changes captured by the target diff:
#
#  1. The single table on the deck (slide 16, "Google Shape;213;p29") is
#     re-styled: its tableStyleId goes from the deck's custom "Table_0"
#     style ({484FA378-C1DB-4A6D-8F2E-E0398C9AA37C}) to the PowerPoint
#     built-in "Medium Style 2" table style
#     ({1C16CF4F-2A5B-4F49-8658-B46C7D23AAFF}).
#
#  2. The presentation's two theme parts have their colour schemes swapped:
#     the deck's live theme (currently "Integral") becomes the stock
#     "Office Theme" palette. (The font scheme / format scheme are already
#     identical between the two themes, so only the 12 scheme colours
#     actually differ.)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------

$targetStyleId  = "{484FA378-C1DB-4A6D-8F2E-E0398C9AA37C}"
$newStyleId     = "{1C16CF4F-2A5B-4F49-8658-B46C7D23AAFF}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable -and $shp.Table.Style -eq $targetStyleId) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme ---------------------------------------------
# RGB() packs as 0x00BBGGRR, matching the PowerPoint object model.
function RGBVal([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),  # 1  dk1
    (RGBVal 0xFF 0xFF 0xFF),  # 2  lt1
    (RGBVal 0x44 0x54 0x6A),  # 3  dk2
    (RGBVal 0xE7 0xE6 0xE6),  # 4  lt2
    (RGBVal 0x5B 0x9B 0xD5),  # 5  accent1
    (RGBVal 0xED 0x7D 0x31),  # 6  accent2
    (RGBVal 0xA5 0xA5 0xA5),  # 7  accent3
    (RGBVal 0xFF 0xC0 0x00),  # 8  accent4
    (RGBVal 0x44 0x72 0xC4),  # 9  accent5
    (RGBVal 0x70 0xAD 0x47),  # 10 accent6
    (RGBVal 0x05 0x63 0xC1),  # 11 hlink
    (RGBVal 0x95 0x4F 0x72)   # 12 folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($c = 1; $c -le $officeThemeColors.Count; $c++) {
    $colorScheme.Colors($c).RGB = $officeThemeColors[$c - 1]
}
